$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update calculated values in row 4
$ws.Range("K4").Value = 1482.95
$ws.Range("L4").Value = 1248.02

# Delete row 10 (After First Recurring duplicate block) and the rest of the
# trailing "Transaction History" sections (rows 12-20), shrinking the used
# range down to A1:T9.
$ws.Range("A10:T20").EntireRow.Delete()
